$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 ---
$ws.Range("A2").Value = "1332119"
$ws.Range("B2").Value = "https://aiesec.org/opportunity/global-talent/1332119"
$ws.Range("C2").Value = "[DSC] Global Sector Development Analyst Intern"
$ws.Range("D2").Value = "Fritz-Erler-Straße 5, 53113 Bonn, Germany"
$ws.Range("E2").Value = "Yes"
$ws.Range("F2").Value = "8 applicants"
$ws.Range("H2").Value = "DHL Group"

# Highlight the PREMIUM cell for this row in yellow (new style/fill)
$ws.Range("E2").Interior.Color = 65535

# --- Row 3 ---
$ws.Range("A3").Value = "1331597"
$ws.Range("B3").Value = "https://aiesec.org/opportunity/global-talent/1331597"
$ws.Range("C3").Value = "Transformation & Value Realisation Intern"
$ws.Range("D3").Value = "Dubai - United Arab Emirates"
$ws.Range("F3").Value = "27 applicants"
$ws.Range("G3").Value = "3 - 6 Months"
$ws.Range("H3").Value = "Dubai Holding Group Services"

# --- Row 4 ---
$ws.Range("A4").Value = "1331591"
$ws.Range("B4").Value = "https://aiesec.org/opportunity/global-talent/1331591"
$ws.Range("C4").Value = "Procurement Intern"
$ws.Range("D4").Value = "Dubai - United Arab Emirates"
$ws.Range("F4").Value = "40 applicants"
$ws.Range("G4").Value = "3 - 6 Months"
$ws.Range("H4").Value = "Dubai Holding Group Services"

# --- Row 5 ---
$ws.Range("A5").Value = "1331590"
$ws.Range("B5").Value = "https://aiesec.org/opportunity/global-talent/1331590"
$ws.Range("C5").Value = "Order to Cash (Accounts Receivable) Intern"
$ws.Range("D5").Value = "Dubai - United Arab Emirates"
$ws.Range("F5").Value = "41 applicants"
$ws.Range("H5").Value = "Dubai Holding Group Services"

# --- Row 6 ---
$ws.Range("A6").Value = "1331468"
$ws.Range("B6").Value = "https://aiesec.org/opportunity/global-talent/1331468"
$ws.Range("C6").Value = "Master Data Management Intern"
$ws.Range("D6").Value = "Dubai - United Arab Emirates"
$ws.Range("F6").Value = "25 applicants"
$ws.Range("G6").Value = "3 - 6 Months"
$ws.Range("H6").Value = "Dubai Holding Group Services"

# --- Row 7 ---
$ws.Range("A7").Value = "1331466"
$ws.Range("B7").Value = "https://aiesec.org/opportunity/global-talent/1331466"
$ws.Range("C7").Value = "Finance Intern"
$ws.Range("D7").Value = "Dubai - United Arab Emirates"
$ws.Range("F7").Value = "16 applicants"
$ws.Range("H7").Value = "Dubai Holding Group Services"

# --- Row 8 ---
$ws.Range("A8").Value = "1331406"
$ws.Range("B8").Value = "https://aiesec.org/opportunity/global-talent/1331406"
$ws.Range("C8").Value = "Front end Developer"
$ws.Range("D8").Value = "Jawhara, Tunisie"
$ws.Range("F8").Value = "2 applicants"
$ws.Range("G8").Value = "9 - 12 Weeks"
$ws.Range("H8").Value = "Digital Glow Agency"

# --- Row 9 (new data) ---
$ws.Range("A9").Value = "1320933"
$ws.Range("B9").Value = "https://aiesec.org/opportunity/global-talent/1320933"
$ws.Range("C9").Value = "Culinary Internship Chef"
$ws.Range("D9").Value = "Hong Kong"
$ws.Range("E9").Value = "No"
$ws.Range("F9").Value = "45 applicants"
$ws.Range("G9").Value = "6 - 18 Months"
$ws.Range("H9").Value = "Treehouse"

# --- Row 10 (new row) ---
$ws.Range("A10").Value = "1289378"
$ws.Range("B10").Value = "https://aiesec.org/opportunity/global-talent/1289378"
$ws.Range("C10").Value = "Medical Advisor (Spanish Speaker)"
$ws.Range("D10").Value = "İstanbul, Türkiye"
$ws.Range("E10").Value = "No"
$ws.Range("F10").Value = "132 applicants"
$ws.Range("G10").Value = "6 - 18 Months"
$ws.Range("H10").Value = "International Plus"

# --- Row 11 (new row) ---
$ws.Range("A11").Value = "1289377"
$ws.Range("B11").Value = "https://aiesec.org/opportunity/global-talent/1289377"
$ws.Range("C11").Value = "Medical Advisor (Italian Speaker)"
$ws.Range("D11").Value = "İstanbul, Türkiye"
$ws.Range("E11").Value = "No"
$ws.Range("F11").Value = "43 applicants"
$ws.Range("G11").Value = "6 - 18 Months"
$ws.Range("H11").Value = "International Plus"

# --- Column width changes ---
# Excel's ColumnWidth setter adds ~5/6 of a character to the stored width
# (internal pixel-padding/MDW rounding), so compensate by subtracting 5/6
# from the desired final width to land exactly on the target value.
$pad = 5 / 6
$ws.Columns.Item(3).ColumnWidth = 49 - $pad
$ws.Columns.Item(4).ColumnWidth = 44 - $pad
$ws.Columns.Item(6).ColumnWidth = 17 - $pad
$ws.Columns.Item(8).ColumnWidth = 31 - $pad
